# Generate Report for Handoff
# Updates the status of the "dd20ad19-..." file from "Handed back: in sync
# with en-US" to "Ready for handoff" across the Overview, zh-cn and de-de
# worksheets, together with refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the dd20ad19-... file ---
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-40-19 04:40:28"

# --- zh-cn sheet: row 3 is the dd20ad19-... file ---
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-19 04:40:25"

# --- de-de sheet: row 3 is the dd20ad19-... file ---
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-19 04:40:28"
